$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 2252.5
$ws.Range("J10").Value = 3005
$ws.Range("L10").Value = 3005
$ws.Range("N10").Value = -3591
$ws.Range("H33").Value = 212.83333
$ws.Range("I33").Value = 125.71429
$ws.Range("J33").Value = 334.8
$ws.Range("K33").Value = 125.71429
$ws.Range("L33").Value = 334.8
$ws.Range("M33").Value = 103.28571
$ws.Range("N33").Value = -792.8
$ws.Range("H51").Value = 2051.8235
$ws.Range("I51").Value = 2480.2
$ws.Range("J51").Value = 1873.3334
$ws.Range("K51").Value = 2480.2
$ws.Range("L51").Value = 1873.3334
$ws.Range("M51").Value = -1996.2
$ws.Range("N51").Value = -2841.3334
$ws.Range("H76").Value = 3779.4119
$ws.Range("J76").Value = 3937.5
$ws.Range("L76").Value = 3937.5
$ws.Range("N76").Value = -4567.5
$ws.Range("H79").Value = 3779.4119
$ws.Range("J79").Value = 3937.5
$ws.Range("L79").Value = 3937.5
$ws.Range("N79").Value = -6121.5
$ws.Range("H98").Value = 2292.3845
$ws.Range("I98").Value = 1556.4286
$ws.Range("J98").Value = 3151
$ws.Range("K98").Value = 1556.4286
$ws.Range("L98").Value = 3151
$ws.Range("M98").Value = -58.42859999999996
$ws.Range("N98").Value = -6147
$ws.Range("H116").Value = 2114.3333
$ws.Range("I116").Value = 2078.6
$ws.Range("J116").Value = 2159
$ws.Range("K116").Value = 2078.6
$ws.Range("L116").Value = 2159
$ws.Range("M116").Value = 1363.4
$ws.Range("N116").Value = -9043
$ws.Range("H122").Value = 2292.3845
$ws.Range("I122").Value = 1556.4286
$ws.Range("J122").Value = 3151
$ws.Range("K122").Value = 4669.2858
$ws.Range("L122").Value = 9453
$ws.Range("M122").Value = -2219.2858
$ws.Range("N122").Value = -14353

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 5000
$ws.Range("K3").Value = 5000
$ws.Range("M3").Value = -4885
$ws.Range("H32").Value = 26722.215
$ws.Range("I32").Value = 26722.215
$ws.Range("K32").Value = 26722.215
$ws.Range("M32").Value = -26435.215
$ws.Range("H74").Value = 2248.4546
$ws.Range("I74").Value = 1895.8334
$ws.Range("J74").Value = 3188.7778
$ws.Range("K74").Value = 1895.8334
$ws.Range("L74").Value = 3188.7778
$ws.Range("M74").Value = -1021.8334
$ws.Range("N74").Value = -4936.7778
$ws.Range("H77").Value = 2248.4546
$ws.Range("I77").Value = 1895.8334
$ws.Range("J77").Value = 3188.7778
$ws.Range("K77").Value = 9479.166999999999
$ws.Range("L77").Value = 15943.889
$ws.Range("M77").Value = -5111.166999999999
$ws.Range("N77").Value = -24679.889
$ws.Range("H132").Value = 6906.875
$ws.Range("I132").Value = 8416.866
$ws.Range("J132").Value = 4390.222
$ws.Range("K132").Value = 25250.598
$ws.Range("L132").Value = 13170.666
$ws.Range("M132").Value = -22720.598
$ws.Range("N132").Value = -18230.666

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4030.5386
$ws.Range("I105").Value = 3989.2222
$ws.Range("J105").Value = 4123.5
$ws.Range("K105").Value = 3989.2222
$ws.Range("L105").Value = 4123.5
$ws.Range("M105").Value = -2242.2222
$ws.Range("N105").Value = -7617.5
$ws.Range("H107").Value = 2921.6667
$ws.Range("I107").Value = 2575
$ws.Range("J107").Value = 3615
$ws.Range("K107").Value = 2575
$ws.Range("L107").Value = 3615
$ws.Range("M107").Value = -655
$ws.Range("N107").Value = -7455

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 6994.5
$ws.Range("I3").Value = 6994.5
$ws.Range("K3").Value = 6994.5
$ws.Range("M3").Value = -6881.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 13929998
$ws.Range("J9").Value = 13929998
$ws.Range("L9").Value = 41789994
$ws.Range("N9").Value = -41790442
$ws.Range("H26").Value = 477.875
$ws.Range("I26").Value = 165.2
$ws.Range("J26").Value = 999
$ws.Range("K26").Value = 495.6
$ws.Range("L26").Value = 2997
$ws.Range("M26").Value = -207.6
$ws.Range("N26").Value = -3573
$ws.Range("H113").Value = 701.60297
$ws.Range("I113").Value = 701.9474
$ws.Range("J113").Value = 699.8182
$ws.Range("K113").Value = 2105.8422
$ws.Range("L113").Value = 2099.4546
$ws.Range("M113").Value = 64.15779999999995
$ws.Range("N113").Value = -6439.4546
$ws.Range("H120").Value = 11611
$ws.Range("I120").Value = 13500
$ws.Range("J120").Value = 7833
$ws.Range("K120").Value = 40500
$ws.Range("L120").Value = 23499
$ws.Range("M120").Value = -35662
$ws.Range("N120").Value = -33175
$ws.Range("H133").Value = 3336.682
$ws.Range("I133").Value = 2805
$ws.Range("J133").Value = 3868.3635
$ws.Range("K133").Value = 8415
$ws.Range("L133").Value = 11605.0905
$ws.Range("M133").Value = -3355
$ws.Range("N133").Value = -21725.0905

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("H41").Value = 9500
$ws.Range("I41").Value = 2000
$ws.Range("K41").Value = 2000
$ws.Range("M41").Value = -1645
$ws.Range("H106").Value = 10100
$ws.Range("J106").Value = 10100
$ws.Range("L106").Value = 10100
$ws.Range("N106").Value = -12624
$ws.Range("H132").Value = 3170.6191
$ws.Range("I132").Value = 2438.0667
$ws.Range("J132").Value = 5002
$ws.Range("K132").Value = 7314.2001
$ws.Range("L132").Value = 15006
$ws.Range("M132").Value = -4784.2001
$ws.Range("N132").Value = -20066
$ws.Range("N4").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3300.4
$ws.Range("I46").Value = 3300.4
$ws.Range("K46").Value = 3300.4
$ws.Range("M46").Value = -3112.4
$ws.Range("H55").Value = 133.5
$ws.Range("J55").Value = 134.5
$ws.Range("L55").Value = 134.5
$ws.Range("N55").Value = -480.5
$ws.Range("H98").Value = 47000
$ws.Range("J98").Value = 47000
$ws.Range("L98").Value = 47000
$ws.Range("N98").Value = -52990

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 20110.334
$ws.Range("I4").Value = 80000
$ws.Range("J4").Value = 2999
$ws.Range("K4").Value = 80000
$ws.Range("L4").Value = 2999
$ws.Range("M4").Value = -79887
$ws.Range("N4").Value = -3225
$ws.Range("H48").Value = 15065
$ws.Range("J48").Value = 15065
$ws.Range("L48").Value = 15065
$ws.Range("N48").Value = -16203
$ws.Range("H93").Value = 67676.5
$ws.Range("I93").Value = 64353
$ws.Range("J93").Value = 71000
$ws.Range("K93").Value = 64353
$ws.Range("L93").Value = 71000
$ws.Range("M93").Value = -61857
$ws.Range("N93").Value = -75992
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("H126").Value = 1661.8462
$ws.Range("I126").Value = 1667
$ws.Range("J126").Value = 1600
$ws.Range("K126").Value = 5001
$ws.Range("L126").Value = 4800
$ws.Range("M126").Value = -2531
$ws.Range("N126").Value = -9740
$ws.Range("N104").ClearContents()
